$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was added as row 351 (data for 2023-01-05, "Segunda"
# quality), which pushes the previous rows 351-375 down to 352-376.
$ws.Rows("351:351").Insert()

$ws.Range("A351").Value = 7
$ws.Range("B351").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C351").Value = "Ñuble"
$ws.Range("D351").Value = 44931
$ws.Range("E351").Value = 16
$ws.Range("F351").Value = 100112023
$ws.Range("G351").Value = "Brócoli"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Segunda"
$ws.Range("J351").Value = 300
$ws.Range("K351").Value = 700
$ws.Range("L351").Value = 700
$ws.Range("M351").Value = 700
$ws.Range("N351").Value = "$/unidad"
$ws.Range("O351").Value = "Región del Maule"
$ws.Range("P351").Value = 700
$ws.Range("Q351").Value = 1
$ws.Range("R351").Value = "Hortaliza"
